$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 1.1.1
$ws.Range("B3").Value = "1.1.1"

# Experimental: (empty) -> "false"
# Plain string assignment of "false" gets auto-typed as a Boolean cell by
# Excel's input parser, so round-trip it through a text formula and then
# "paste values" to collapse it back down to a literal shared-string cell
# (keeps the original cell style too).
$cell = $ws.Range("B7")
$cell.Formula = "=""false"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Date: 2022-05-04T10:16:52-05:00 -> 2022-10-21T09:04:31-05:00
$ws.Range("B8").Value = "2022-10-21T09:04:31-05:00"
